$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add note about pool quality column (PoolQC header is in G13)
$ws.Range("H13").Value = "泳池质量，删除/重建"

# Add note about Fireplaces column (header is in D27)
$ws.Range("E27").Value = "缺失值太多删除"

# Update the active cell selection on Sheet1
$ws.Activate()
$ws.Range("E24").Select()
